# Add two new student-mark rows (5 and 6) to the worksheet, matching the
# layout of the existing data rows (USN + 9 groups of Internals/Externals/
# Total/Remarks columns spanning C:AL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is the last existing data row; copy its formatting (the bordered /
# bold / centered style used on column A) down into the two new rows so the
# new rows look consistent with the rest of the table.
$ws.Range("A4:AL4").Copy()
$ws.Range("A5:AL5").PasteSpecial(-4122)
$ws.Range("A6:AL6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row5 = @(3, "1AM18CS010", 26, 29, 55, "P", 28, 21, 49, "P", 23, 35, 58, "P", 34, 15, 49, "F", 33, 18, 51, "F", 21, 33, 54, "P", 39, 27, 66, "P", 25, 41, 66, "P", 37, 29, 66, "P")
$row6 = @(4, "1AM18CS028", 20, 21, 41, "P", 21, 0, 21, "A", 22, 8, 30, "F", 18, 15, 33, "F", 19, 9, 28, "F", 17, 9, 26, "F", 33, 27, 60, "P", 25, 38, 63, "P", 37, 26, 63, "P")

for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
    $ws.Cells.Item(6, $i + 1).Value = $row6[$i]
}
